$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "OrSplitWithLoop" - copy of "OrSplitWithAndSplit", appended at the end.
#    (created first so it gets the lower sheetId, matching the target file)
# ---------------------------------------------------------------------------
$srcLoop = $wb.Worksheets.Item("OrSplitWithAndSplit")
$srcLoop.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$loop = $wb.Worksheets.Item($wb.Worksheets.Count)
$loop.Name = "OrSplitWithLoop"

# Trim the AndSplit branch (rows 8,10,11,12,13) down to a single Loop branch.
# Delete from the bottom up so row numbers of earlier deletes stay valid.
$loop.Rows.Item(13).Delete()
$loop.Rows.Item(12).Delete()
$loop.Rows.Item(11).Delete()
$loop.Rows.Item(10).Delete()
$loop.Rows.Item(8).Delete()

$loop.Range("A7").Value = "Loop"
$loop.Range("B8").Value = "TestItem_OrSplit:0"
$loop.Range("C8").Value = "LoopRight"

$loop.Range("C8").Select()

# ---------------------------------------------------------------------------
# 2. "SequenceWithOrSplit" - copy of "StartWithOrSplit", inserted right after
#    "StartWithOrSplit" (position 3).
# ---------------------------------------------------------------------------
$after = $wb.Worksheets.Item("StartWithOrSplit")
$srcSeq = $wb.Worksheets.Item("StartWithOrSplit")
$srcSeq.Copy($null, $after)
$seq = $wb.Worksheets.Item(3)
$seq.Name = "SequenceWithOrSplit"

# Rework the existing 10 rows into the inner OrSplit block first (this is
# what introduces the new shared strings "Left1"/"Left2"), then turn row 3
# into the wrapping "First" activity, then append the trailing
# Right2/End/End/Last rows.
$seq.Range("A4").Value = "OrSplit"
$seq.Range("A5").Value = "Block"
$seq.Range("B5").Value = ""
$seq.Range("C5").Value = ""
$seq.Range("A6").Value = "Elementary"
$seq.Range("B6").Value = "TestItem_OrSplit:0"
$seq.Range("C6").Value = "Left1"
$seq.Range("A7").Value = "Elementary"
$seq.Range("B7").Value = "TestItem_OrSplit:0"
$seq.Range("C7").Value = "Left2"
$seq.Range("A8").Value = "BlockEnd"
$seq.Range("B8").Value = ""
$seq.Range("C8").Value = ""
$seq.Range("A9").Value = "Block"
$seq.Range("A10").Value = "Elementary"
$seq.Range("B10").Value = "TestItem_OrSplit:0"
$seq.Range("C10").Value = "Right1"

$seq.Range("A3").Value = "Elementary"
$seq.Range("B3").Value = "TestItem_OrSplit:0"
$seq.Range("C3").Value = "First"

$seq.Range("A10:C10").Copy()
$seq.Range("A11:C14").PasteSpecial(-4122)

$seq.Range("A11").Value = "Elementary"
$seq.Range("B11").Value = "TestItem_OrSplit:0"
$seq.Range("C11").Value = "Right2"
$seq.Range("A12").Value = "End"
$seq.Range("B12").Value = ""
$seq.Range("C12").Value = ""
$seq.Range("A13").Value = "End"
$seq.Range("B13").Value = ""
$seq.Range("C13").Value = ""
$seq.Range("A14").Value = "Elementary"
$seq.Range("B14").Value = "TestItem_OrSplit:0"
$seq.Range("C14").Value = "Last"

$seq.Range("C14").Select()
$seq.Select()

Write-Host "done"
